$d = $word.ActiveDocument

# --- Change 1: insert a new sub-bullet before
#     "Does molecular diet give you different kinds of links (e.g. functional groups) than other approaches?"
#     The new bullet sits at the deeper list level (ilvl=1 / ListLevelNumber=2),
#     matching the sibling bullets around it (numId=1). ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Does molecular diet give you different kinds of links")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($targetIndex)
    $newPara.Range.Text = "Proportion of total that they are interacting with (corrected by species richness)"
    # Demote from the level it inherited (ilvl=0) down to ilvl=1 to match the diff.
    $newPara.Range.ListFormat.ListIndent()
}

# --- Change 2: mark a rendered page break at the start of the run that begins
#     "Sample cleaning, denoising, taxonomic assignment, and rarefying methods:" ---
$cleanIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Sample cleaning, denoising, taxonomic assignment, and rarefying methods:")) {
        $cleanIndex = $i
        break
    }
}

if ($cleanIndex -ge 1) {
    $cleanPara = $d.Paragraphs.Item($cleanIndex)
    $paraId = $cleanPara.ParaId
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' " +
           "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='$paraId'>" +
           "<w:r><w:lastRenderedPageBreak/>" +
           "<w:t>Sample cleaning, denoising, taxonomic assignment, and rarefying methods:</w:t>" +
           "</w:r></w:p>"
    $cleanPara.Range.InsertXML($xml)
}
